# GlueX computing model: update first 10 rows based on recent (Oct. 2016) estimates.
#
# 1) Event-reconstruction CPU time per event (row 9) now derived from David's
#    Oct. 2016 benchmark (250Hz/16 threads Ivy Bridge, 340Hz/24 threads
#    Haswell, 600Hz/36 threads Broadwell) instead of the old fixed 1/22.
# 2) Raw event size (row 11) is now computed from the Spring 2016 data fit
#    (16.1/11.5 kB base + 4.6/2.3 kB per 10^7 g/s) instead of a flat 18000
#    bytes assumption.
# 3) Updated the explanatory comment text (column F) for both rows.
# 4) Moved the active-cell selection on the "model" sheet to A15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# --- Row 9: CPU time per event (CPU-s/event) -------------------------------
$ws.Range("B9").Formula = "=1/(340/24)"
$ws.Range("C9").Formula = "=1/(340/24)"
$ws.Range("D9").Formula = "=1/(600/36)"
$ws.Range("F9").Value = "Oct. 2016 David benchmark gave 250Hz for 16 threads on Ivy Bridge, 340Hz for 24 threads on Haswell, 600Hz for 36 threads on Broadwell"

# --- Row 11: raw event size (bytes) -----------------------------------------
$ws.Range("B11").Formula = "=(16.1+4.6*1)*1000"
$ws.Range("C11").Formula = "=(11.5+2.3*2)*1000"
$ws.Range("D11").Formula = "=(11.5+2.3*5)*1000"
$ws.Range("F11").Value = "size of a single raw event. Actual Spring 2016 data is 16.4kB+4.6kB/10^7 g/s. Estimate from Spring 2016 data for reduced windows is 11.5kB + 0.23kB/10^7 g/s"

# --- Move the on-sheet selection to A15 -------------------------------------
$ws.Activate()
$ws.Range("A15").Select()

# --- Best-effort: restore the saved window position (xWindow) --------------
$excel.ActiveWindow.Left = 600

$wb.Application.Calculate()
